$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row (row 1)
$ws.Range("E1").Value = "symbol"
$ws.Range("F1").Value = "name"

# Update data row (row 2)
$ws.Range("B2").Value = "check that user can search by new added GF>> Name"
$ws.Range("C2").Value = "METLVNLIVASFLYKLGLFSSLGVSQSHYVKANGLSTTTKLSSICKTSDLTIHKKSNRTRKFSVSAGYRDGSRSGSSGDFIAGFLLGGAVFGAVAYIFAPQIRRSVLNEEDEYGFEKPKQPTYYDEGLEKTRETLNEKIGQLNSAIDNVSSRLRGREKNTSSLNVPVETDPEVEATT"
$ws.Range("D2").Value = "Gene"
$ws.Range("E2").Value = "GF_name_test"
$ws.Range("F2").Value = "name_GF"

# Column width change for column E (OOXML width target 17.85546875;
# ColumnWidth maps to OOXML width as width = ColumnWidth + 5/6, rounded to
# the nearest 1/6 character, so use the closest achievable value)
$ws.Range("E1").ColumnWidth = 17.0221354166667

# Selection / view changes
$ws.Range("C9").Select()
